$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.053, 0.414, 1.124, -33, 37, 55),
    @(-0.119, 0.09, 0.089, -90, -68, 12),
    @(-0.201, 0.099, -0.644, -25, -6, 126),
    @(-0.521, -0.095, 0.976, 98, 1, 52),
    @(0.145, 0.163, -0.84, 94, 45, -175),
    @(0.027, -0.168, -0.511, 163, -13, -122),
    @(0.157, 0.668, -0.052, -84, 31, -95),
    @(-0.954, 0.062, -0.051, 59, 51, -161),
    @(0.404, 0.388, -0.567, 0, 17, -173),
    @(0.459, -0.616, 0.576, -89, 11, 129)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}
